$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 3.3
$ws.Range("I2").Value = 2.28
$ws.Range("J2").Value = 3.85
$ws.Range("L2").Value = 1.36
$ws.Range("M2").Value = 1.05
$ws.Range("N2").Value = 4.5
$ws.Range("O2").Value = 1.26
$ws.Range("P2").Value = 2.18
$ws.Range("Q2").Value = 1.79
$ws.Range("U2").Value = 2.38
$ws.Range("V2").Value = 1.78
$ws.Range("X2").Value = 18.5
$ws.Range("AB2").Value = 15
$ws.Range("AF2").Value = 25
$ws.Range("AN2").Value = 28
# Row 3
$ws.Range("F3").Value = 2.96
$ws.Range("G3").Value = 3.4
$ws.Range("H3").Value = 2.32
$ws.Range("I3").Value = 2.6
$ws.Range("J3").Value = 3.4
$ws.Range("K3").Value = 3.9
$ws.Range("M3").Value = 1.06
$ws.Range("N3").Value = 3.45
$ws.Range("O3").Value = 1.31
$ws.Range("P3").Value = 1.95
$ws.Range("Q3").Value = 1.86
$ws.Range("R3").Value = 1.37
$ws.Range("S3").Value = 2.84
$ws.Range("U3").Value = 2.12
$ws.Range("V3").Value = 1.62
$ws.Range("W3").Value = 1.42
$ws.Range("Y3").Value = 13
$ws.Range("AC3").Value = 10.5
# Row 4
$ws.Range("P4").Value = 2
$ws.Range("Q4").Value = 1.83
# Row 5
$ws.Range("H5").Value = 14
$ws.Range("I5").Value = 17.5
$ws.Range("L5").Value = 1.41
$ws.Range("P5").Value = 1.8
$ws.Range("Q5").Value = 2
$ws.Range("U5").Value = 1.56
$ws.Range("Y5").Value = 40
$ws.Range("AB5").Value = 6.6
$ws.Range("AC5").Value = 15
$ws.Range("AD5").Value = 1000
$ws.Range("AF5").Value = 6.8
$ws.Range("AG5").Value = 13.5
$ws.Range("AJ5").Value = 10.5
$ws.Range("AK5").Value = 19
$ws.Range("AN5").Value = 7.8
# Row 6
$ws.Range("F6").Value = 1.95
$ws.Range("G6").Value = 2.02
$ws.Range("W6").Value = 1.98
# Row 7
$ws.Range("F7").Value = 1.4
$ws.Range("G7").Value = 1.48
$ws.Range("H7").Value = 7.2
$ws.Range("J7").Value = 5.2
$ws.Range("L7").Value = 1.21
$ws.Range("P7").Value = 2.76
$ws.Range("Q7").Value = 1.47
$ws.Range("T7").Value = 1.7
$ws.Range("W7").Value = 3.05
$ws.Range("Y7").Value = 990
$ws.Range("AD7").Value = 32
# Row 8
$ws.Range("F8").Value = 1.62
$ws.Range("G8").Value = 1.72
$ws.Range("K8").Value = 3.75
$ws.Range("Q8").Value = 2.7
$ws.Range("V8").Value = 1.13
$ws.Range("AF8").Value = 8.4
$ws.Range("AG8").Value = 12
$ws.Range("AH8").Value = 1000
$ws.Range("AN8").Value = 980
# Row 9
$ws.Range("F9").Value = 1.9
$ws.Range("G9").Value = 1.93
$ws.Range("J9").Value = 3.6
$ws.Range("K9").Value = 3.85
$ws.Range("P9").Value = 1.74
$ws.Range("V9").Value = 1.24
$ws.Range("W9").Value = 2.06
# Row 10
$ws.Range("F10").Value = 1.8
$ws.Range("G10").Value = 1.87
$ws.Range("H10").Value = 5
$ws.Range("J10").Value = 3.55
$ws.Range("K10").Value = 3.9
$ws.Range("P10").Value = 1.86
$ws.Range("Q10").Value = 1.83
$ws.Range("T10").Value = 1.87
$ws.Range("V10").Value = 1.22
$ws.Range("W10").Value = 2.12
$ws.Range("X10").Value = 1000
# Row 11
$ws.Range("F11").Value = 2.16
$ws.Range("G11").Value = 2.34
$ws.Range("H11").Value = 3.4
$ws.Range("I11").Value = 3.8
$ws.Range("J11").Value = 3.35
$ws.Range("L11").Value = 1.41
$ws.Range("N11").Value = 3.4
$ws.Range("O11").Value = 1.33
$ws.Range("P11").Value = 1.85
$ws.Range("Q11").Value = 1.98
$ws.Range("R11").Value = 1.31
$ws.Range("S11").Value = 3.5
$ws.Range("T11").Value = 1.81
$ws.Range("U11").Value = 2.04
$ws.Range("V11").Value = 1.36
$ws.Range("W11").Value = 1.74
$ws.Range("X11").Value = 16
$ws.Range("Z11").Value = 980
$ws.Range("AA11").Value = 75
$ws.Range("AB11").Value = 10
$ws.Range("AD11").Value = 16
$ws.Range("AE11").Value = 50
$ws.Range("AH11").Value = 19
$ws.Range("AI11").Value = 60
$ws.Range("AK11").Value = 980
$ws.Range("AN11").Value = 20
$ws.Range("AO11").Value = 1000
# Row 12
$ws.Range("F12").Value = 2.38
$ws.Range("H12").Value = 3.3
$ws.Range("I12").Value = 3.6
$ws.Range("J12").Value = 3.25
$ws.Range("K12").Value = 3.4
$ws.Range("N12").Value = 3.2
$ws.Range("P12").Value = 1.75
$ws.Range("T12").Value = 1.86
$ws.Range("W12").Value = 1.67
# Row 13
$ws.Range("F13").Value = 1.95
$ws.Range("G13").Value = 2.02
$ws.Range("I13").Value = 4.7
$ws.Range("P13").Value = 1.78
$ws.Range("Q13").Value = 2.06
$ws.Range("U13").Value = 1.92
$ws.Range("W13").Value = 1.98
$ws.Range("AC13").Value = 1000
